$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "set_reference"
$ws.Range("B1").Value = "order"

$ws.Range("A2").Value = "HJC1000007197"
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = "GENT1000019888"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "GENT1000006796"
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = "HJCN1002228"
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = "BOU1000004753"
$ws.Range("B6").Value = 5

$ws.Range("A1:B6").Select()
